$d = $word.ActiveDocument
foreach ($f in $d.Fields) {
  $r = $f.Code
  # "self" starts after " m:" (3 chars) within this 13-char range " m:self.name "
  $txt = $r.Text
  Write-Output ("text=[" + $txt + "] len=" + $txt.Length)
  $idx = $txt.IndexOf("self")
  Write-Output ("idx=" + $idx)
  $subStart = $r.Start + $idx
  $subEnd = $subStart + 4
  $sub = $d.Range($subStart, $subEnd)
  Write-Output ("sub text=[" + $sub.Text + "]")
  Write-Output ("color before=" + $sub.Font.Color)
  $sub.Font.Color = $sub.Font.Color
  Write-Output ("color after=" + $sub.Font.Color)
}
